# Round every numeric metric value in the results table (rows 4-12,
# columns B:AE) to 2 decimal places, matching the "round 2" pass from
# the commit. Columns that were already at <=2 decimals (K,L,M,Q,T,W,X,Y)
# are unaffected because rounding them is a no-op.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 4
$lastRow = 12
$firstCol = 2   # column B
$lastCol = 31   # column AE

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $raw = $cell.Value2
        if ($raw -ne $null) {
            $cell.Value = [Math]::Round([double]$raw, 2)
        }
    }
}
